$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 528
$ws1.Range("F9").Value = 286
$ws1.Range("F10").Value = 3075

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 528
$ws4.Range("F10").Value = 286
$ws4.Range("F11").Value = 3075
